$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header / label strings (values unchanged, only text relabeled)
$ws.Range("B1").Value = "e_modulus"
$ws.Range("C1").Value = "tensile_yield_strength"
$ws.Range("D1").Value = "tensile_strain_at_break"

$ws.Range("A2").Value = "e_modulus"
$ws.Range("A3").Value = "tensile_yield_strength"
$ws.Range("A4").Value = "tensile_strain_at_break"

# Update the active selection to match the saved view state
$ws.Range("B10").Select()
